# "WIP: building nicer rundown"
# Renames several header cells on the TeamA / TeamB / Rundown sheets to
# shorter, English labels, shrinks the now-narrower "# A" / "# B" columns
# on the Rundown sheet, and updates the remembered cell selection on a
# few sheets (cosmetic state that Excel persists on save).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TeamA
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("TeamA")
$wsA.Activate()
$wsA.Range("A1").Value = "#"
$wsA.Range("C1").Value = "Played"
$wsA.Range("C2").Select()

# ---------------------------------------------------------------------
# TeamB
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("TeamB")
$wsB.Activate()
$wsB.Range("A1").Value = "#"
$wsB.Range("C1").Value = "Played"
$wsB.Range("C2").Select()

# ---------------------------------------------------------------------
# Rundown
# ---------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Rundown")
$wsR.Activate()
$wsR.Range("A1").Value = "# A"
$wsR.Range("B1").Value = "Score A"
$wsR.Range("D1").Value = "# B"
$wsR.Range("E1").Value = "Score B"

# The "# A"/"# B" headers (and the numbers underneath them) are much
# narrower than the old "A Nummer"/"B Nummer" labels, so the best-fit
# columns shrink.
$wsR.Columns.Item(1).ColumnWidth = 3
$wsR.Columns.Item(4).ColumnWidth = 3

$wsR.Range("C1").Select()

# Keep Rundown the active sheet/tab, as it was before the edit.
$wsR.Activate()
